# Reading locators from excel file now. Added one demo test to run.
# Populate the "OR" locator sheet (rows 2-9) with the News / FrontEnd
# locator definitions, then update the active selection and the
# width of column B to match the new (narrower) content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OR")

# Columns: A=PageName, B=Locator, C=ObjectType, D=Property (xpath)
$rows = @(
    @{ Row = 2; PageName = "News";          Locator = "NewsSubmissionPage.NewsTitle";        ObjectType = "Textbox";  Property = '//label[contains(text(),"Title")]/following-sibling::input' },
    @{ Row = 3; PageName = "News";          Locator = "NewsSubmissionPage.NewsSummary";      ObjectType = "Textbox";  Property = "//b[contains(text(),'Summary')]//..//..//textarea[@class='pt-input pt-fixed']" },
    @{ Row = 4; PageName = "News";          Locator = "NewsSubmissionPage.NewsBody";         ObjectType = "Textarea"; Property = "//label[contains(.,`"Body`")]//following::div[@class='id-froala']//div[@contenteditable='true']" },
    @{ Row = 5; PageName = "News";          Locator = "NewsSubmissionPage.NewsSubmitArticle"; ObjectType = "Button";  Property = '//div[@class=''pt-fill'' and contains(text(),"Submit Article")]' },
    @{ Row = 6; PageName = "News";          Locator = "NewsSubmissionPage.NewsSubmit";       ObjectType = "Button";   Property = "//div[@class='id-form__footer']/button[@type='button']" },
    @{ Row = 7; PageName = "FrontEndLogin"; Locator = "FrontEndLoginButton";                 ObjectType = "Button";   Property = '//*[@id="btnLogin"]' },
    @{ Row = 8; PageName = "FrontEndLogin"; Locator = "FrontEndPasswordTextBox";              ObjectType = "Textbox";  Property = '//*[@id="textPassword"]' },
    @{ Row = 9; PageName = "FrontEndLogin"; Locator = "FrontEndUsernameTextBox";              ObjectType = "Textbox";  Property = '//*[@id="textUsername"]' }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.PageName
    $ws.Cells.Item($r.Row, 2).Value = $r.Locator
    $ws.Cells.Item($r.Row, 3).Value = $r.ObjectType
    $ws.Cells.Item($r.Row, 4).Value = $r.Property
}

# D8 ends up with the default/normal style (its shading was cleared while
# editing), unlike the other data cells in this block.
$ws.Range("D8").Style = "Normal"

# Narrow column B now that it holds shorter "Locator" names instead of
# long descriptive text.
$ws.Columns.Item(2).ColumnWidth = 37.3

# Move/restore the active selection to A9, as recorded in the workbook
# view state.
$ws.Range("A9").Select()
